# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# New columns AD:AF are appended right after the existing "Unnamed: 28" (AC)
# column: headers in row 1, then a constant W/L/T for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
# Grab the formatting of the last existing header cell (AC1, which carries
# the bold/bordered/centered header style) and apply it to the three new
# header cells before setting their text, so we don't clobber the labels
# with the copied cell's own value.
$headerFmtSrc = $ws.Range("AC1")
$newHeaders = $ws.Range("AD1:AF1")
$headerFmtSrc.Copy()
$newHeaders.PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-47) ------------------------------------------------
# Every player on the roster shares the team's overall record for the
# season, so the same three numbers repeat down each column.
$ws.Range("AD2:AD47").Value = 86
$ws.Range("AE2:AE47").Value = 76
$ws.Range("AF2:AF47").Value = 0
